$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the explicit "justify" (w:jc w:val="both") paragraph
#    alignment override from every paragraph in the main body. The
#    justification now comes from the (updated) Normal style instead,
#    so setting the COM alignment to the Word default (left) causes the
#    explicit <w:jc/> to disappear from each paragraph's <w:pPr/> -- and
#    the whole (now-empty) <w:pPr/> disappears too when jc was its only
#    child, matching the target markup exactly.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $p.Range.ParagraphFormat.Alignment = 0
}

# ---------------------------------------------------------------------
# 2) Tiny text tweak in the "Conexión ..." bullet: the run split around
#    the (hidden) _GoBack bookmark moves two characters to the right --
#    " per" + "mitirá..." becomes " permi" + "tirá...". The rendered
#    text is identical ("...Esto permitirá que..."), only where the
#    runs/bookmark divide changes.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$pos = $bm.Start

# Insert the two extra characters right at the bookmark, then delete
# the now-duplicated "mi" that used to start the following run.
$insPoint = $d.Range($pos, $pos)
$insPoint.InsertBefore("mi")

# Re-separate the newly typed "mi" (plus the " per" text before it)
# from its neighbours so saving doesn't silently re-merge the runs --
# a harmless bold on/off toggle is enough to force Word to keep the
# run boundary.
$permi = $d.Range($pos - 4, $pos + 2)
$permi.Font.Bold = 1
$permi.Font.Bold = 0

# Remove the stale "mi" that originally led the following run (now
# shifted two characters to the right by the insert above).
$stale = $d.Range($pos + 2, $pos + 4)
$stale.Text = ""

# ---------------------------------------------------------------------
# 3) Give the "Normal" style an explicit 709-twip (0.5 cm / 35.45 pt)
#    first-line indent plus justified alignment, so paragraphs based on
#    it (i.e. all of them) inherit the justification that used to be
#    hard-coded per-paragraph.
# ---------------------------------------------------------------------
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.FirstLineIndent = 35.45
$normal.ParagraphFormat.Alignment = 3
